$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.886.05"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "1.626.73"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.72"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.69%  "
$ws.Range("E9").Value = "  +3.47%  "
$ws.Range("E10").Value = "  +2.44%  "
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "1.859.37"
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("D13").Value = "1.629.91"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("E14").Value = "  +6.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +23.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.18%  "
$ws.Range("D17").Value = "29.891.95"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.82%  "
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.29%  "
$ws.Range("E27").Value = "  +2.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +2.89%  "
$ws.Range("E31").Value = "  +6.20%  "
$ws.Range("E32").Value = "  +4.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("D34").Value = "1.429.47"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  +7.37%  "
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("E39").Value = "  +3.11%  "
$ws.Range("E40").Value = "  +3.07%  "
$ws.Range("B41").Value = "BitcoinSV"
$ws.Range("C41").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "55.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.830"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.93%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.98%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0494"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("D49").Value = "1.767.29"
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "89.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.56%  "
$ws.Range("D51").Value = "0.0₆0108"
$ws.Range("E51").Value = "  +5.07%  "
